$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings (e.g. "19.40", "0.0847")
# must be forced to Text format first, otherwise Excel COM auto-converts them to
# numbers and mangles the exact text (dropping trailing zeros, using scientific notation, etc).
$numericTextCells = @(
    "D5"
    "D10"
    "D11"
    "D16"
    "D20"
    "D25"
    "D27"
    "D29"
    "D39"
    "D40"
    "D41"
    "D46"
    "D47"
)
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply all cell value updates as described by the diff.
$ws.Range("D2").Value = "26.332.89"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "1.592.77"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "211.69"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "19.40"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "1.817.07"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "1.586.38"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "64.59"
$ws.Range("D17").Value = "26.346.91"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("E19").Value = "  +3.08%  "
$ws.Range("D20").Value = "212.16"
$ws.Range("E20").Value = "  +2.81%  "
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("D25").Value = "144.86"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").Value = "7.06"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("D29").Value = "15.21"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("D34").Value = "1.337.64"
$ws.Range("E34").Value = "  +4.29%  "
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "1.04"
$ws.Range("E39").Value = "  -17.39%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "0.818"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("D41").Value = "5.79"
$ws.Range("E41").Value = "  +5.07%  "
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").Value = "1.729.31"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").Value = "61.70"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").Value = "87.88"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("E50").Value = "  -2.87%  "
$ws.Range("E51").Value = "  -0.70%  "

# Restore default ("Normal") style on the cells we temporarily switched to Text format,
# so no stray style index is introduced for those cells.
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).Style = "Normal"
}
